# Auto-generated Excel COM-interop script
# Applies numeric corrections to the LeveProfit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H19").Value = 2535.818
$ws.Range("I19").Value = 2097.75
$ws.Range("J19").Value = 2786.1428
$ws.Range("K19").Value = 2097.75
$ws.Range("L19").Value = 2786.1428
$ws.Range("M19").Value = -1922.75
$ws.Range("N19").Value = -3136.1428

$ws.Range("H55").Value = 250.85185
$ws.Range("I55").Value = 256.86365
$ws.Range("J55").Value = 224.4
$ws.Range("K55").Value = 256.86365
$ws.Range("L55").Value = 224.4
$ws.Range("M55").Value = -42.86365000000001
$ws.Range("N55").Value = -652.4

$ws.Range("H106").Value = 53332.668
$ws.Range("I106").Value = 63199.6
$ws.Range("K106").Value = 63199.6
$ws.Range("M106").Value = -62568.6


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H5").Value = 780.3333
$ws.Range("J5").Value = 950
$ws.Range("L5").Value = 950
$ws.Range("N5").Value = -1174

$ws.Range("H74").Value = 201629.9
$ws.Range("I74").Value = 316528.9
$ws.Range("K74").Value = 316528.9
$ws.Range("M74").Value = -315654.9

$ws.Range("H77").Value = 201629.9
$ws.Range("I77").Value = 316528.9
$ws.Range("K77").Value = 1582644.5
$ws.Range("M77").Value = -1578276.5

$ws.Range("H97").Value = 1543888.8
$ws.Range("I97").Value = 1950079.4
$ws.Range("J97").Value = 364.8
$ws.Range("K97").Value = 1950079.4
$ws.Range("L97").Value = 364.8
$ws.Range("M97").Value = -1949583.4
$ws.Range("N97").Value = -1356.8

$ws.Range("H132").Value = 3585.2
$ws.Range("I132").Value = 3150.5454
$ws.Range("J132").Value = 10757
$ws.Range("K132").Value = 9451.636200000001
$ws.Range("L132").Value = 32271
$ws.Range("M132").Value = -6921.636200000001
$ws.Range("N132").Value = -37331


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H4").Value = 780.3333
$ws.Range("J4").Value = 950
$ws.Range("L4").Value = 950
$ws.Range("N4").Value = -1180

$ws.Range("H80").Value = 52095.95
$ws.Range("I80").Value = 413.6
$ws.Range("J80").Value = 69323.39999999999
$ws.Range("K80").Value = 413.6
$ws.Range("L80").Value = 69323.39999999999
$ws.Range("M80").Value = 584.4
$ws.Range("N80").Value = -71319.39999999999

$ws.Range("H83").Value = 52095.95
$ws.Range("I83").Value = 413.6
$ws.Range("J83").Value = 69323.39999999999
$ws.Range("K83").Value = 2068
$ws.Range("L83").Value = 346617
$ws.Range("M83").Value = 2924
$ws.Range("N83").Value = -356601

$ws.Range("H105").Value = 33343416
$ws.Range("I105").Value = 45467532
$ws.Range("J105").Value = 2102.875
$ws.Range("K105").Value = 45467532
$ws.Range("L105").Value = 2102.875
$ws.Range("M105").Value = -45465785
$ws.Range("N105").Value = -5596.875

$ws.Range("H123").Value = 77499.5
$ws.Range("J123").Value = 77499.5
$ws.Range("L123").Value = 77499.5
$ws.Range("N123").Value = -87299.5


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 321.2143
$ws.Range("I7").Value = 230.53847
$ws.Range("K7").Value = 230.53847
$ws.Range("M7").Value = -117.53847

$ws.Range("H58").Value = 3144.5625
$ws.Range("I58").Value = 1474.65
$ws.Range("J58").Value = 5927.75
$ws.Range("K58").Value = 1474.65
$ws.Range("L58").Value = 5927.75
$ws.Range("M58").Value = -1271.65
$ws.Range("N58").Value = -6333.75

$ws.Range("H64").Value = 45000
$ws.Range("I64").Value = 20000
$ws.Range("K64").Value = 20000
$ws.Range("M64").Value = -19752

$ws.Range("H67").Value = 45000
$ws.Range("I67").Value = 20000
$ws.Range("K67").Value = 20000
$ws.Range("M67").Value = -19142

$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80630

$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -82184

$ws.Range("H105").Value = 10513.4
$ws.Range("I105").Value = 21178
$ws.Range("J105").Value = 3403.6667
$ws.Range("K105").Value = 21178
$ws.Range("L105").Value = 3403.6667
$ws.Range("M105").Value = -19431
$ws.Range("N105").Value = -6897.6667

$ws.Range("H107").Value = 831.9091
$ws.Range("I107").Value = 706.6667
$ws.Range("J107").Value = 982.2
$ws.Range("K107").Value = 706.6667
$ws.Range("L107").Value = 982.2
$ws.Range("M107").Value = 1213.3333
$ws.Range("N107").Value = -4822.2

$ws.Range("H134").Value = 4266.9355
$ws.Range("I134").Value = 3381.1428
$ws.Range("J134").Value = 6127.1
$ws.Range("K134").Value = 10143.4284
$ws.Range("L134").Value = 18381.3
$ws.Range("M134").Value = -7608.428400000001
$ws.Range("N134").Value = -23451.3

$ws.Range("H136").Value = 3144.5625
$ws.Range("I136").Value = 1474.65
$ws.Range("J136").Value = 5927.75
$ws.Range("K136").Value = 4423.950000000001
$ws.Range("L136").Value = 17783.25
$ws.Range("M136").Value = -1873.950000000001
$ws.Range("N136").Value = -22883.25


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H9").Value = 5648.8
$ws.Range("J9").Value = 6623.5
$ws.Range("L9").Value = 19870.5
$ws.Range("N9").Value = -20318.5

$ws.Range("H11").Value = 5355.174
$ws.Range("J11").Value = 2871.5
$ws.Range("L11").Value = 8614.5
$ws.Range("N11").Value = -8894.5

$ws.Range("H12").Value = 175.19048
$ws.Range("I12").Value = 100.75
$ws.Range("K12").Value = 302.25
$ws.Range("M12").Value = -129.25

$ws.Range("H34").Value = 1858.9565
$ws.Range("J34").Value = 2120.6
$ws.Range("L34").Value = 6361.799999999999
$ws.Range("N34").Value = -6529.799999999999

$ws.Range("H57").Value = 1199.6666
$ws.Range("J57").Value = 1199.5
$ws.Range("L57").Value = 3598.5
$ws.Range("N57").Value = -4716.5

$ws.Range("H61").Value = 59.35294
$ws.Range("J61").Value = 64.29031999999999
$ws.Range("L61").Value = 192.87096
$ws.Range("N61").Value = -622.87096

$ws.Range("H116").Value = 1513.5
$ws.Range("I116").Value = 1513.5
$ws.Range("K116").Value = 4540.5
$ws.Range("M116").Value = -1098.5

$ws.Range("H117").Value = 1952.1428
$ws.Range("I117").Value = 219.33333
$ws.Range("J117").Value = 3251.75
$ws.Range("K117").Value = 657.99999
$ws.Range("L117").Value = 9755.25
$ws.Range("M117").Value = 2784.00001
$ws.Range("N117").Value = -16639.25

$ws.Range("H132").Value = 3480.0454
$ws.Range("J132").Value = 4126.364
$ws.Range("L132").Value = 37137.276
$ws.Range("N132").Value = -42197.276


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H2").Value = 675
$ws.Range("I2").Value = 715.625
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 715.625
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -602.625
$ws.Range("N2").Value = -576

$ws.Range("H69").Value = 22200
$ws.Range("J69").Value = 22200
$ws.Range("L69").Value = 22200
$ws.Range("N69").Value = -23698

$ws.Range("H70").Value = 8548.299999999999
$ws.Range("I70").Value = 4110.75
$ws.Range("J70").Value = 11506.667
$ws.Range("K70").Value = 4110.75
$ws.Range("L70").Value = 11506.667
$ws.Range("M70").Value = -3840.75
$ws.Range("N70").Value = -12046.667

$ws.Range("H72").Value = 22200
$ws.Range("J72").Value = 22200
$ws.Range("L72").Value = 66600
$ws.Range("N72").Value = -74088

$ws.Range("H73").Value = 8548.299999999999
$ws.Range("I73").Value = 4110.75
$ws.Range("J73").Value = 11506.667
$ws.Range("K73").Value = 4110.75
$ws.Range("L73").Value = 11506.667
$ws.Range("M73").Value = -3174.75
$ws.Range("N73").Value = -13378.667

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H113").Value = 19893.5
$ws.Range("J113").Value = 2126.5
$ws.Range("L113").Value = 2126.5
$ws.Range("N113").Value = -6466.5

$ws.Range("H132").Value = 16500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 16500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 49500
$ws.Range("N132").Value = -54560
$ws.Range("M132").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 5283.8486
$ws.Range("I46").Value = 1427.3
$ws.Range("J46").Value = 6960.609
$ws.Range("K46").Value = 1427.3
$ws.Range("L46").Value = 6960.609
$ws.Range("M46").Value = -1239.3
$ws.Range("N46").Value = -7336.609

$ws.Range("H132").Value = 6546.0527
$ws.Range("I132").Value = 4469.1177
$ws.Range("J132").Value = 24200
$ws.Range("K132").Value = 13407.3531
$ws.Range("L132").Value = 72600
$ws.Range("M132").Value = -10877.3531
$ws.Range("N132").Value = -77660

$ws.Range("H136").Value = 3763.4412
$ws.Range("I136").Value = 2570.6296
$ws.Range("K136").Value = 7711.888800000001
$ws.Range("M136").Value = -5161.888800000001


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H8").Value = 2166
$ws.Range("I8").Value = 498
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 498
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = -358
$ws.Range("N8").Value = -3280

$ws.Range("H31").Value = 35000
$ws.Range("I31").Value = 36666.668
$ws.Range("J31").Value = 30000
$ws.Range("K31").Value = 36666.668
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = -36318.668
$ws.Range("N31").Value = -30696

